# Apply crypto price/volume updates per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NumberFormat "@" (Text) is set before each write so that numeric-looking
# strings (e.g. "0.0980", "5.10", "55.182.31") are stored as literal text,
# matching the original inlineStr cell type instead of being auto-coerced
# into floating point numbers by Excel.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '55.182.31'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.94%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.286.60'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.73%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '505.89'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +1.53%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '129.33'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.995'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.46%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.94%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.310.60'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +1.59%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0980'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +3.15%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +1.35%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.10'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +8.50%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.342'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +2.43%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '23.64'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +4.57%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.697.45'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.94%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '55.279.84'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +2.16%  '
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +1.59%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.296.67'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.76%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.44'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +2.33%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.19'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +1.56%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '312.67'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +3.25%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.61'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +4.75%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.997'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.31%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '60.22'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -1.27%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.993'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.98%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +4.01%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.53'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +3.28%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '172.23'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.93%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0₃0712'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +4.48%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.14'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +4.01%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.15'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +7.21%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +1.83%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '18.01'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +1.69%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.50%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +4.18%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.912'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -3.98%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +5.72%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +2.36%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.46'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +4.69%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.376'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +1.07%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '136.53'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +9.34%  '
$ws.Range('B43').NumberFormat = '@'
$ws.Range('B43').Value = 'Filecoin'
$ws.Range('C43').NumberFormat = '@'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.49'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +4.01%  '
$ws.Range('B44').NumberFormat = '@'
$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').NumberFormat = '@'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.92'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +2.35%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '260.95'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +9.49%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0507'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +3.05%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.551'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +1.43%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.375'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +1.18%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +3.88%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '16.59'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +2.89%  '
